# Deploy 28d4dc4 from branch master
#
# Summary of edits applied to public/AddInfos/en/3.5.1.xlsx:
#  - Tabelle1 (sheet 1):
#      * B10 text "224 - 226" -> "224 - 266" (fixes a typo, matches the
#        actual source-data max value of 266 on Tabelle2!D6)
#      * A34/B34 (previously blank) now read "Specification:" / "Estimated
#        data", styled like the other footnote-label cells in col A/B
#      * A35 text "Datenquelle:" -> "Source:" (English wording)
#      * Selection moved from J16 to K30
#  - Tabelle2 (sheet 2): selection moved from D1:D7 to D8
#  - Chart1: category-axis title "Jahr" -> "Year"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Tabelle1 cell content updates -----------------------------------
# New footnote row: "Specification:" / "Estimated data"
$ws1.Range("A34").Value = "Specification:"
$ws1.Range("B34").Value = "Estimated data"
# Match the font used by the other label cells in that footer block
# (A35:C36 etc. use the 10pt font / style index 7).
$ws1.Range("A34:B34").Font.Size = 10

# "Datenquelle:" -> "Source:"
$ws1.Range("A35").Value = "Source:"

# Fix the mislabeled range for year 2015 to match the real data (224-266)
$ws1.Range("B10").Value = "224 - 266"

# --- Chart axis title translation ------------------------------------
$co = $ws1.ChartObjects().Item(1)
$chart = $co.Chart()
$catAx = $chart.Axes(1)
$catAx.AxisTitle().Text = "Year"

# --- Selection / view state -------------------------------------------
$ws2.Range("D8").Select()
$ws1.Activate()
$ws1.Range("K30").Select()
